# Updated capital structure database
# Applies the diff: updates row 2 and row 3 values for the first two companies,
# and removes row 4 (Arbuthnot Banking Group PLC) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (first company) updates
# ---------------------------------------------------------------------------

# B2: change the numeric-looking label from "2" to "1", but keep it stored as
# text (matching the original inlineStr type) rather than letting Excel
# auto-convert it to a number.
$c = $ws.Range("B2")
$c.NumberFormat = "@"
$c.Value2 = "1"
$c.Style = "Normal"

# D2 is removed entirely in the new data.
$ws.Range("D2").ClearContents()

$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -524.2
$ws.Range("L2").Value = -1.73003300330033
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0.0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0.0

# T2 is removed entirely.
$ws.Range("T2").ClearContents()

$ws.Range("U2").Value = 3813.6
$ws.Range("V2").Value = 11.57038834951456
$ws.Range("W2").Value = -0.2334550636857576
$ws.Range("X2").Value = 0.1396363601256733
$ws.Range("Y2").Value = -0.3730914238114309
$ws.Range("Z2").Value = 0.1921345322189953
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04654910020482791
$ws.Range("AC2").Value = -0.04654910020482791
$ws.Range("AD2").Value = 1438.8
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1438.8
$ws.Range("AG2").Value = -2374.8
$ws.Range("AH2").Value = 0.813616828771771
$ws.Range("AI2").Value = 0.463142985900985
$ws.Range("AJ2").Value = 1.16115783297477
$ws.Range("AK2").Value = 3.358981612446958

# AN2 and AP2 are removed entirely.
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# ---------------------------------------------------------------------------
# Row 3 (second company) updates
# ---------------------------------------------------------------------------

$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -524.2
$ws.Range("L3").Value = -1.73003300330033
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0

$ws.Range("U3").Value = 3813.6
$ws.Range("V3").Value = 11.57038834951456
$ws.Range("W3").Value = -0.2334550636857576
$ws.Range("X3").Value = 0.1396363601256733
$ws.Range("Y3").Value = -0.3730914238114309
$ws.Range("Z3").Value = 0.1921345322189953
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04654910020482791
$ws.Range("AC3").Value = -0.04654910020482791
$ws.Range("AD3").Value = 1438.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1438.8
$ws.Range("AG3").Value = -2374.8
$ws.Range("AH3").Value = 0.813616828771771
$ws.Range("AI3").Value = 0.463142985900985
$ws.Range("AJ3").Value = 1.16115783297477
$ws.Range("AK3").Value = 3.358981612446958

# AN3 and AP3 are removed entirely.
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# ---------------------------------------------------------------------------
# Row 4 (Arbuthnot Banking Group PLC) is removed entirely; the dimension
# shrinks from A1:AQ4 to A1:AQ3 as a result.
# ---------------------------------------------------------------------------
$ws.Rows("4").Delete()
